$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a rolling weekly price log. Each week, 3 new records (for the
# latest reporting date) are inserted at the top (after the header row), and
# every existing data row shifts down by 3. The weekly update does not trim
# the tail, so the 3 oldest rows that fall off the "front window" are simply
# appended, unchanged, at the very end of the sheet.

$firstDataRow = 2
$lastDataRow  = 56
$numCols      = 18   # columns A..R
$shift        = 3

# 1) Capture the current data rows (2..56) across columns A..R before
#    overwriting anything.
$data = @()
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $row = @()
    for ($c = 1; $c -le $numCols; $c++) {
        $row += ,($ws.Cells.Item($r, $c).Value2)
    }
    $data += ,$row
}

$dateNumberFormat = $ws.Cells.Item($firstDataRow, 4).NumberFormat

function Set-DataRow($rowIndex, $values) {
    for ($c = 1; $c -le $numCols; $c++) {
        $ws.Cells.Item($rowIndex, $c).Value2 = $values[$c - 1]
    }
    $ws.Cells.Item($rowIndex, 4).NumberFormat = $dateNumberFormat
}

# 2) Shift every captured row down by 3 (process from the bottom up so we
#    never clobber a source row before it has been read/copied).
for ($i = ($data.Count - 1); $i -ge 0; $i--) {
    $srcRow = $firstDataRow + $i
    $dstRow = $srcRow + $shift
    Set-DataRow $dstRow $data[$i]
}

# 3) Write the 3 brand-new weekly records into the freed-up rows 2..4.
Set-DataRow 2 @(10, "Vega Modelo de Temuco", "La Araucanía", 44515, 9, 300000000, "Espárragos", "Sin especificar", "Extra",   150, 1500, 1500, 1500, "`$/kilo", "Región del Maule", 1500, 1, "Hortaliza")
Set-DataRow 3 @(10, "Vega Modelo de Temuco", "La Araucanía", 44515, 9, 300000000, "Espárragos", "Sin especificar", "Primera", 300, 1300, 1300, 1300, "`$/kilo", "Región del Maule", 1300, 1, "Hortaliza")
Set-DataRow 4 @(10, "Vega Modelo de Temuco", "La Araucanía", 44515, 9, 300000000, "Espárragos", "Sin especificar", "Segunda", 400, 1000, 1000, 1000, "`$/kilo", "Región del Maule", 1000, 1, "Hortaliza")

Write-Host "Applied weekly roll: 3 new rows inserted, data shifted, sheet now spans to row" ($lastDataRow + $shift)
